$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (component U1 x ...): quantity 10 -> 1, cost total 2 -> 0.2
$ws.Range("A9").Value = 1
$ws.Range("J9").Value = 0.2

# Row 10: quantity 10 -> 1, cost total 2 -> 0.2
$ws.Range("A10").Value = 1
$ws.Range("J10").Value = 0.2

# Row 11: quantity 10 -> 1, cost total 2 -> 0.2
$ws.Range("A11").Value = 1
$ws.Range("J11").Value = 0.2

# Row 12: quantity 10 -> 1, cost total 2 -> 0.2
$ws.Range("A12").Value = 1
$ws.Range("J12").Value = 0.2

# J19 is =SUM(J7:J18); it recalculates automatically to 29.669999999999998

# Move the active selection to K12, matching the author's last selection
$ws.Range("K12").Select()
